# Scheduled-runner style refresh of cached market-price / profit columns
# (H:currentAveragePrice, I:currentAveragePriceNQ, J:currentAveragePriceHQ,
#  K:LevePriceNQ, L:LevePriceHQ, M:LeveProfitNQ, N:LeveProfitHQ)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1507.0625
$ws.Range("I40").Value = 1337.7142
$ws.Range("K40").Value = 1337.7142
$ws.Range("M40").Value = -1162.7142

$ws.Range("H116").Value = 9632.637000000001
$ws.Range("I116").Value = 10808.375
$ws.Range("J116").Value = 6497.3335
$ws.Range("K116").Value = 10808.375
$ws.Range("L116").Value = 6497.3335
$ws.Range("M116").Value = -7366.375
$ws.Range("N116").Value = -13381.3335

$ws.Range("H138").Value = 4793.6562
$ws.Range("I138").Value = 1767.5454
$ws.Range("J138").Value = 6378.7617
$ws.Range("K138").Value = 5302.6362
$ws.Range("L138").Value = 19136.2851
$ws.Range("M138").Value = -162.6361999999999
$ws.Range("N138").Value = -29416.2851

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4728.8945
$ws.Range("I32").Value = 4354.2285
$ws.Range("K32").Value = 4354.2285
$ws.Range("M32").Value = -4067.2285

$ws.Range("H45").Value = 9529.375
$ws.Range("I45").Value = 16711.143
$ws.Range("K45").Value = 16711.143
$ws.Range("M45").Value = -16334.143

$ws.Range("H61").Value = 3030.9062
$ws.Range("I61").Value = 2818.077
$ws.Range("K61").Value = 2818.077
$ws.Range("M61").Value = -2606.077

$ws.Range("H122").Value = 1744.7333
$ws.Range("I122").Value = 1356.6364
$ws.Range("K122").Value = 4069.9092
$ws.Range("M122").Value = -1619.9092

$ws.Range("H124").Value = 42900
$ws.Range("J124").Value = 42900
$ws.Range("L124").Value = 42900
$ws.Range("N124").Value = -52720

$ws.Range("H132").Value = 4286.524
$ws.Range("I132").Value = 3724.125
$ws.Range("K132").Value = 11172.375
$ws.Range("M132").Value = -8642.375

$ws.Range("H136").Value = 3030.9062
$ws.Range("I136").Value = 2818.077
$ws.Range("K136").Value = 8454.231
$ws.Range("M136").Value = -5904.231

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 16431.334
$ws.Range("J81").Value = 16431.334
$ws.Range("L81").Value = 16431.334
$ws.Range("N81").Value = -18553.334

$ws.Range("H84").Value = 16431.334
$ws.Range("J84").Value = 16431.334
$ws.Range("L84").Value = 49294.00199999999
$ws.Range("N84").Value = -59902.00199999999

$ws.Range("H94").Value = 1838.4814
$ws.Range("I94").Value = 1539.5416
$ws.Range("K94").Value = 1539.5416
$ws.Range("M94").Value = -1088.5416

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 950
$ws.Range("J16").Value = 1000
$ws.Range("L16").Value = 1000
$ws.Range("N16").Value = -1574

$ws.Range("H58").Value = 1905.931
$ws.Range("I58").Value = 1917.72
$ws.Range("K58").Value = 1917.72
$ws.Range("M58").Value = -1714.72

$ws.Range("H99").Value = 3817.6667

$ws.Range("H113").Value = 950
$ws.Range("J113").Value = 1000
$ws.Range("L113").Value = 1000
$ws.Range("N113").Value = -5340

$ws.Range("H126").Value = 3817.6667

$ws.Range("H132").Value = 4755.077
$ws.Range("I132").Value = 4386.5
$ws.Range("K132").Value = 13159.5
$ws.Range("M132").Value = -10629.5

$ws.Range("H134").Value = 15845.6
$ws.Range("I134").Value = 11934.5
$ws.Range("K134").Value = 35803.5
$ws.Range("M134").Value = -33268.5

$ws.Range("H135").Value = 74257.14
$ws.Range("J135").Value = 74257.14
$ws.Range("L135").Value = 74257.14
$ws.Range("N135").Value = -84397.14

$ws.Range("H136").Value = 1905.931
$ws.Range("I136").Value = 1917.72
$ws.Range("K136").Value = 5753.16
$ws.Range("M136").Value = -3203.16

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9099.799999999999
$ws.Range("J3").Value = 12833.167
$ws.Range("L3").Value = 38499.501
$ws.Range("N3").Value = -38723.501

$ws.Range("H80").Value = 4056.5
$ws.Range("J80").Value = 5368.6665
$ws.Range("L80").Value = 16105.9995
$ws.Range("N80").Value = -17977.9995

$ws.Range("H83").Value = 4056.5
$ws.Range("J83").Value = 5368.6665
$ws.Range("L83").Value = 48317.9985
$ws.Range("N83").Value = -57677.9985

$ws.Range("H86").Value = 634.75
$ws.Range("I86").Value = 600
$ws.Range("J86").Value = 669.5
$ws.Range("K86").Value = 1800
$ws.Range("L86").Value = 2008.5
$ws.Range("M86").Value = -614
$ws.Range("N86").Value = -4380.5

$ws.Range("H89").Value = 634.75
$ws.Range("I89").Value = 600
$ws.Range("J89").Value = 669.5
$ws.Range("K89").Value = 5400
$ws.Range("L89").Value = 6025.5
$ws.Range("M89").Value = 528
$ws.Range("N89").Value = -17881.5

$ws.Range("H107").Value = 954.6667
$ws.Range("I107").Value = 265
$ws.Range("K107").Value = 795
$ws.Range("M107").Value = 1125

$ws.Range("H113").Value = 396.0476
$ws.Range("J113").Value = 363.5625
$ws.Range("L113").Value = 1090.6875
$ws.Range("N113").Value = -5430.6875

$ws.Range("H120").Value = 15000
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 15000
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 45000
$ws.Range("M120").ClearContents()
$ws.Range("N120").Value = -54676

$ws.Range("H133").Value = 16142.714
$ws.Range("J133").Value = 12428.571
$ws.Range("L133").Value = 37285.713
$ws.Range("N133").Value = -47405.713

$ws.Range("H141").Value = 117923
$ws.Range("I141").Value = 7659.25
$ws.Range("K141").Value = 22977.75
$ws.Range("M141").Value = -17797.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 717782.4
$ws.Range("I132").Value = 912451.25
$ws.Range("J132").Value = 3996.3333
$ws.Range("K132").Value = 2737353.75
$ws.Range("L132").Value = 11988.9999
$ws.Range("M132").Value = -2734823.75
$ws.Range("N132").Value = -17048.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 10009
$ws.Range("J29").Value = 10009
$ws.Range("L29").Value = 10009
$ws.Range("N29").Value = -10599

$ws.Range("H40").Value = 4651.846
$ws.Range("I40").Value = 4246.625
$ws.Range("J40").Value = 5300.2
$ws.Range("K40").Value = 4246.625
$ws.Range("L40").Value = 5300.2
$ws.Range("M40").Value = -4110.625
$ws.Range("N40").Value = -5572.2

$ws.Range("H61").Value = 3808.1052
$ws.Range("I61").Value = 3696.8125
$ws.Range("J61").Value = 4401.6665
$ws.Range("K61").Value = 3696.8125
$ws.Range("L61").Value = 4401.6665
$ws.Range("M61").Value = -3494.8125
$ws.Range("N61").Value = -4805.6665

$ws.Range("H113").Value = 3808.1052
$ws.Range("I113").Value = 3696.8125
$ws.Range("J113").Value = 4401.6665
$ws.Range("K113").Value = 3696.8125
$ws.Range("L113").Value = 4401.6665
$ws.Range("M113").Value = -1526.8125
$ws.Range("N113").Value = -8741.666499999999

$ws.Range("H135").Value = 124749.5
$ws.Range("J135").Value = 124749.5
$ws.Range("L135").Value = 124749.5
$ws.Range("N135").Value = -134889.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 731.1539
$ws.Range("J113").Value = 798.5
$ws.Range("L113").Value = 2395.5
$ws.Range("N113").Value = -6735.5

$ws.Range("H131").Value = 96416.336
$ws.Range("J131").Value = 96416.336
$ws.Range("L131").Value = 96416.336
$ws.Range("N131").Value = -106496.336

$ws.Range("H132").Value = 3755.276
$ws.Range("I132").Value = 3762.24
$ws.Range("K132").Value = 11286.72
$ws.Range("M132").Value = -8642.375
